$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2024-03-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-30 Saturday", 2)

# Update each division-problem table cell by row/column position
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "43÷4="
$t.Cell(1, 2).Range.Text = "93÷7="
$t.Cell(1, 3).Range.Text = "97÷3="
$t.Cell(1, 4).Range.Text = "36÷5="
$t.Cell(1, 5).Range.Text = "64÷7="
$t.Cell(5, 1).Range.Text = "27÷3="
$t.Cell(5, 2).Range.Text = "64÷3="
$t.Cell(5, 3).Range.Text = "44÷6="
$t.Cell(5, 4).Range.Text = "27÷4="
$t.Cell(5, 5).Range.Text = "59÷9="
$t.Cell(9, 1).Range.Text = "41÷3="
$t.Cell(9, 2).Range.Text = "35÷5="
$t.Cell(9, 3).Range.Text = "21÷9="
$t.Cell(9, 4).Range.Text = "58÷3="
$t.Cell(9, 5).Range.Text = "62÷2="
$t.Cell(13, 1).Range.Text = "92÷2="
$t.Cell(13, 2).Range.Text = "78÷2="
$t.Cell(13, 3).Range.Text = "58÷7="
$t.Cell(13, 4).Range.Text = "19÷3="
$t.Cell(13, 5).Range.Text = "42÷5="
$t.Cell(17, 1).Range.Text = "13÷3="
$t.Cell(17, 2).Range.Text = "37÷9="
$t.Cell(17, 3).Range.Text = "19÷4="
$t.Cell(17, 4).Range.Text = "76÷7="
$t.Cell(17, 5).Range.Text = "46÷9="
